$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.595.00"
$ws.Range("E2").Value = "  +2.18%  "

$ws.Range("D3").Value = "'2.361.04"
$ws.Range("E3").Value = "  +1.90%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'522.63"
$ws.Range("E5").Value = "  +0.68%  "

$ws.Range("D6").Value = "'136.74"
$ws.Range("E6").Value = "  +1.99%  "

$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").Value = "'0.540"
$ws.Range("E8").Value = "  +0.59%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.103"
$ws.Range("E9").Value = "  -0.40%  "

$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'5.46"
$ws.Range("E10").Value = "  +5.32%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.153"
$ws.Range("E11").Value = "  -1.23%  "

$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "'0.344"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'24.39"
$ws.Range("E13").Value = "  +1.26%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'2.781.07"
$ws.Range("E14").Value = "  +1.80%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "'57.554.41"
$ws.Range("E15").Value = "  +1.93%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000136"
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'2.361.70"
$ws.Range("E17").Value = "  +2.18%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'10.66"
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'330.56"
$ws.Range("E19").Value = "  +2.88%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.26"
$ws.Range("E20").Value = "  -0.92%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.74"
$ws.Range("E21").Value = "  +1.44%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'61.62"
$ws.Range("E23").Value = "  +1.23%  "

$ws.Range("B24").Value = "Kaspa"
$ws.Range("C24").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D24").Value = "'0.166"
$ws.Range("E24").Value = "  +4.08%  "

$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "  +0.46%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'8.32"
$ws.Range("E26").Value = "  +8.25%  "

$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").Value = "'1.33"
$ws.Range("E27").Value = "  +9.36%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "'0.0₃0748"
$ws.Range("E28").Value = "  +1.56%  "

$ws.Range("D29").Value = "'170.41"
$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.71"
$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "'6.30"
$ws.Range("E31").Value = "  +0.42%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'18.62"
$ws.Range("E32").Value = "  +0.83%  "

$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'1.31"
$ws.Range("E34").Value = "  +3.50%  "

$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "'0.993"
$ws.Range("E35").Value = "  +0.24%  "

$ws.Range("B36").Value = "SuiNetwork"
$ws.Range("C36").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D36").Value = "'0.932"
$ws.Range("E36").Value = "  -0.81%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.06"
$ws.Range("E37").Value = "  +0.76%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'1.59"
$ws.Range("E38").Value = "  +4.30%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'38.68"
$ws.Range("E39").Value = "  +3.19%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'151.88"
$ws.Range("E40").Value = "  +7.64%  "

$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.383"
$ws.Range("E41").Value = "  +0.29%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'3.67"
$ws.Range("E42").Value = "  +1.92%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'5.36"
$ws.Range("E43").Value = "  +4.77%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "'284.23"
$ws.Range("E44").Value = "  +0.78%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.0941"
$ws.Range("E45").Value = "  +1.20%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0510"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.567"
$ws.Range("E47").Value = "  +1.50%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'18.43"
$ws.Range("E48").Value = "  +6.52%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0222"
$ws.Range("E49").Value = "  +2.25%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'17.75"
$ws.Range("E50").Value = "  +4.37%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'1.67"
$ws.Range("E51").Value = "  +0.42%  "
